$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7, columns B..O (col A, the rank index, is unchanged)
$data = @(
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(1,1,1,1,1,1,1,1,1,1,1,1,1,1),
    @(5,3,5,5,5,5,3,5,3,3,5,3,3,3),
    @(3,4,3,3,3,3,5,3,2,2,3,5,5,2),
    @(2,5,2,4,2,4,2,2,5,5,2,2,2,5),
    @(4,2,4,2,4,2,4,4,4,4,4,4,4,4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colIndex = $j + 2
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}
